$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 8
$ws.Range("C9").Value = 21

$ws.Range("E44").Value = "Yes"
$ws.Range("E45").Value = "Yes"
$ws.Range("E46").Value = "Yes"
$ws.Range("E48").Value = "Yes"
$ws.Range("E49").Value = "Yes"
$ws.Range("E50").Value = "Yes"

$ws.Range("C8").Select()
